$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.015279619371116
$ws.Range("D2").Value = 1.021040710356193
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.013604256292397
$ws.Range("I2").Value = 1.025974687004851
$ws.Range("J2").Value = 1.020506041757966
$ws.Range("K2").Value = 1.023879214411177
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.016464909508728
$ws.Range("N2").Value = 1.021955277390941

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.016287228064566
$ws.Range("D3").Value = 1.021757400558965
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.015256553053639
$ws.Range("I3").Value = 1.02613572338924
$ws.Range("J3").Value = 1.021148288838729
$ws.Range("K3").Value = 1.024402518411296
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.017919614021175
$ws.Range("N3").Value = 1.022598436536226

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.016938298029711
$ws.Range("D4").Value = 1.022219932620392
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.016324587441439
$ws.Range("I4").Value = 1.026237568325745
$ws.Range("J4").Value = 1.021562406150504
$ws.Range("K4").Value = 1.024739252578872
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.018859349185116
$ws.Range("N4").Value = 1.023013141942083

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.017211789047437
$ws.Range("D5").Value = 1.022414090510346
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.01677333081261
$ws.Range("I5").Value = 1.026279819722075
$ws.Range("J5").Value = 1.021736152632825
$ws.Range("K5").Value = 1.024880366265855
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.01925404825016
$ws.Range("N5").Value = 1.023187135164339

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.017257696605926
$ws.Range("D6").Value = 1.022446673437919
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.016848661883474
$ws.Range("I6").Value = 1.026286880828537
$ws.Range("J6").Value = 1.02176530502709
$ws.Range("K6").Value = 1.024904033534405
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.01932029872064
$ws.Range("N6").Value = 1.023216328958347

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.016941953287956
$ws.Range("D7").Value = 1.022222528109318
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.0163305845774
$ws.Range("I7").Value = 1.026238135108109
$ws.Range("J7").Value = 1.02156472912746
$ws.Range("K7").Value = 1.024741139912138
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.018864624601399
$ws.Range("N7").Value = 1.023015468217933

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.015620336904306
$ws.Range("D8").Value = 1.021283170716676
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.014162892788172
$ws.Range("I8").Value = 1.026029597855483
$ws.Range("J8").Value = 1.020723394955062
$ws.Range("K8").Value = 1.024056456632299
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.016956860037399
$ws.Range("N8").Value = 1.022172939254508

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.013284362901566
$ws.Range("D9").Value = 1.019618593493405
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.010334267894711
$ws.Range("I9").Value = 1.025644084639166
$ws.Range("J9").Value = 1.019229633780589
$ws.Range("K9").Value = 1.022835552181105
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.013582908378075
$ws.Range("N9").Value = 1.02067705676783

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.0117221609921
$ws.Range("D10").Value = 1.018502599671376
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.007775381769413
$ws.Range("I10").Value = 1.025374945575503
$ws.Range("J10").Value = 1.018226181113466
$ws.Range("K10").Value = 1.022011901804862
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.011324950350499
$ws.Range("N10").Value = 1.019672179082823

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.011044527780292
$ws.Range("D11").Value = 1.018017868420655
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.00666570882549
$ws.Range("I11").Value = 1.025255527711981
$ws.Range("J11").Value = 1.017789853780064
$ws.Range("K11").Value = 1.021652940400705
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.010345078425815
$ws.Range("N11").Value = 1.01923523211456

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.010792643755041
$ws.Range("D12").Value = 1.017837591922212
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.006253267956569
$ws.Range("I12").Value = 1.025210737929175
$ws.Range("J12").Value = 1.017627506607879
$ws.Range("K12").Value = 1.021519257503279
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.009980777141811
$ws.Range("N12").Value = 1.019072654390769

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.010846681953858
$ws.Range("D13").Value = 1.017876272067556
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.006341749797939
$ws.Range("I13").Value = 1.025220365071547
$ws.Range("J13").Value = 1.017662343130792
$ws.Range("K13").Value = 1.021547948728697
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.01005893622053
$ws.Range("N13").Value = 1.019107540385541

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.01102371067074
$ws.Range("D14").Value = 1.018002971317369
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.006631621673526
$ws.Range("I14").Value = 1.025251834201871
$ws.Range("J14").Value = 1.017776439746519
$ws.Range("K14").Value = 1.021641897247723
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.010314972015242
$ws.Range("N14").Value = 1.019221799031549

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.011132759872174
$ws.Range("D15").Value = 1.018081004977385
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.006810186698544
$ws.Range("I15").Value = 1.025271166033671
$ws.Range("J15").Value = 1.017846701849748
$ws.Range("K15").Value = 1.021699735807404
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.010472679707145
$ws.Range("N15").Value = 1.019292160915024

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.011767108030913
$ws.Range("D16").Value = 1.018534738065638
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.007848991335323
$ws.Range("I16").Value = 1.025382810278583
$ws.Range("J16").Value = 1.01825510013021
$ws.Range("K16").Value = 1.022035676042071
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.011389934962818
$ws.Range("N16").Value = 1.019701139167888

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.012164697889622
$ws.Range("D17").Value = 1.018818951231194
$ws.Range("E17").Value = 0.9894763578477731
$ws.Range("F17").Value = 1.008500155167841
$ws.Range("I17").Value = 1.025452070904792
$ws.Range("J17").Value = 1.018510787750605
$ws.Range("K17").Value = 1.022245782037465
$ws.Range("L17").Value = 0.9930127773692701
$ws.Range("M17").Value = 1.011964720338567
$ws.Range("N17").Value = 1.019957189894031

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.012396490709975
$ws.Range("D18").Value = 1.018984583415254
$ws.Range("E18").Value = 0.9897087662937551
$ws.Range("F18").Value = 1.00887980867141
$ws.Range("I18").Value = 1.025492191698699
$ws.Range("J18").Value = 1.01865974996191
$ws.Range("K18").Value = 1.022368109967625
$ws.Range("L18").Value = 0.9932001317071766
$ws.Range("M18").Value = 1.012299774979522
$ws.Range("N18").Value = 1.020106363648763

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.012475506695347
$ws.Range("D19").Value = 1.019041035208884
$ws.Range("E19").Value = 0.9897880325774039
$ws.Range("F19").Value = 1.009009233972017
$ws.Range("I19").Value = 1.025505824744233
$ws.Range("J19").Value = 1.01871051240904
$ws.Range("K19").Value = 1.022409782771057
$ws.Range("L19").Value = 0.993264023964098
$ws.Range("M19").Value = 1.012413984939175
$ws.Range("N19").Value = 1.02015719818439

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.01212205213272
$ws.Range("D20").Value = 1.018788472804418
$ws.Range("E20").Value = 0.9894336180360677
$ws.Range("F20").Value = 1.008430308005636
$ws.Range("I20").Value = 1.025444668626285
$ws.Range("J20").Value = 1.018483373110057
$ws.Range("K20").Value = 1.022223262760285
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.011903072879706
$ws.Range("N20").Value = 1.019929736321549

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.010971585136685
$ws.Range("D21").Value = 1.01796566778229
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.006546268857666
$ws.Range("I21").Value = 1.025242579272493
$ws.Range("J21").Value = 1.017742848761687
$ws.Range("K21").Value = 1.0216142413814
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.010239585122615
$ws.Range("N21").Value = 1.019188160343665

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.010247193267747
$ws.Range("D22").Value = 1.017447031295701
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.005360195684187
$ws.Range("I22").Value = 1.025113014284694
$ws.Range("J22").Value = 1.017275655421572
$ws.Range("K22").Value = 1.021229307879892
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.009191752151536
$ws.Range("N22").Value = 1.018720303535419

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.010631307138359
$ws.Range("D23").Value = 1.017722094337535
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.005989101267055
$ws.Range("I23").Value = 1.025181936484263
$ws.Range("J23").Value = 1.017523475262638
$ws.Range("K23").Value = 1.021433559928354
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.009747414251442
$ws.Range("N23").Value = 1.018968475309085

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.012141322268741
$ws.Range("D24").Value = 1.018802245140287
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.008461869423214
$ws.Range("I24").Value = 1.025448014255595
$ws.Range("J24").Value = 1.018495761150032
$ws.Range("K24").Value = 1.022233438942453
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.011930929354579
$ws.Range("N24").Value = 1.019942141953962

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.013889121857044
$ws.Range("D25").Value = 1.020050031416563
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.011325162436934
$ws.Range("I25").Value = 1.025745887129323
$ws.Range("J25").Value = 1.019617143543852
$ws.Range("K25").Value = 1.023152895198623
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.014456649401944
$ws.Range("N25").Value = 1.021065116839406
